$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value2 -eq "6-28-2012-13") {
        # Force the new value to be stored as literal text - otherwise Excel's
        # smart-entry parser would reinterpret "2013-06-28" as a date serial.
        $cell.NumberFormat = "@"
        $cell.Value = "2013-06-28"
    }
}
